# Weekly fruit/vegetable price update: a new weekly record is inserted
# above the old row 300 (date serial 44753 = 2022-07-11), pushing the
# existing rows 300-315 down to 301-316.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 300; everything below shifts down by one.
$ws.Rows.Item(300).Insert()

# Populate the newly inserted row 300 with the new record's data.
$ws.Range("A300").Value = 3
$ws.Range("B300").Value = "Femacal de La Calera"
$ws.Range("C300").Value = "Coquimbo"
$ws.Range("D300").Value = 44753
$ws.Range("E300").Value = 5
$ws.Range("F300").Value = 100112001
$ws.Range("G300").Value = "Berenjena"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 105
$ws.Range("K300").Value = 8500
$ws.Range("L300").Value = 9000
$ws.Range("M300").Value = 8738
$ws.Range("N300").Value = "`$/caja 60 unidades"
$ws.Range("O300").Value = "Región de Arica y Parinacota"
$ws.Range("P300").Value = 146
$ws.Range("Q300").Value = 60
$ws.Range("R300").Value = "Hortaliza"
